$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds text values (coming from a scraper, not real
# numbers) even when they look numeric, e.g. "212.86". Excel's COM Value
# setter auto-detects such strings and silently turns the cell into a
# Number, which would lose the original text typing (and things like
# trailing zeros). Force those specific cells to Text format first so the
# new value is stored as a string, matching the source data.
$textCells = @("D5","D8","D11","D16","D18","D22","D25","D31","D32","D40","D43","D45","D48","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.938.21"
$ws.Range("E2").Value = "  +1.41%  "

$ws.Range("D3").Value = "1.642.24"
$ws.Range("E3").Value = "  +0.80%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").Value = "212.86"
$ws.Range("E5").Value = "  +0.76%  "

$ws.Range("E6").Value = "  +0.61%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").Value = "23.46"
$ws.Range("E8").Value = "  +2.17%  "

$ws.Range("E9").Value = "  -1.25%  "

$ws.Range("E10").Value = "  +0.76%  "

$ws.Range("D11").Value = "0.0883"
$ws.Range("E11").Value = "  +2.33%  "

$ws.Range("D12").Value = "1.875.19"
$ws.Range("E12").Value = "  +0.86%  "

$ws.Range("D13").Value = "1.641.61"

$ws.Range("E14").Value = "  +1.25%  "

$ws.Range("E15").Value = "  +2.69%  "

$ws.Range("D16").Value = "65.58"
$ws.Range("E16").Value = "  +0.81%  "

$ws.Range("D17").Value = "27.938.39"
$ws.Range("E17").Value = "  +1.50%  "

$ws.Range("D18").Value = "233.22"
$ws.Range("E18").Value = "  +1.92%  "

$ws.Range("D19").Value = "0.0₃0723"
$ws.Range("E19").Value = "  +0.70%  "

$ws.Range("E20").Value = "  +0.99%  "

$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").Value = "10.43"
$ws.Range("E22").Value = "  -2.91%  "

$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("E24").Value = "  -1.83%  "

$ws.Range("D25").Value = "153.15"
$ws.Range("E25").Value = "  +2.88%  "

$ws.Range("E26").Value = "  +0.72%  "

$ws.Range("E27").Value = "  +0.82%  "

$ws.Range("E28").Value = "  +0.17%  "

$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("E30").Value = "  +0.81%  "

$ws.Range("D31").Value = "0.0484"
$ws.Range("E31").Value = "  +0.88%  "

$ws.Range("D32").Value = "3.38"
$ws.Range("E32").Value = "  +3.53%  "

$ws.Range("E33").Value = "  +0.56%  "

$ws.Range("D34").Value = "1.407.47"
$ws.Range("E34").Value = "  -3.77%  "

$ws.Range("E35").Value = "  +2.51%  "

$ws.Range("E36").Value = "  +1.59%  "

$ws.Range("E37").Value = "  +1.90%  "

$ws.Range("E38").Value = "  +0.65%  "

$ws.Range("E39").Value = "  +0.78%  "

$ws.Range("D40").Value = "0.926"
$ws.Range("E40").Value = "  +0.85%  "

$ws.Range("E41").Value = "  +0.85%  "

$ws.Range("E42").Value = "  +0.05%  "

$ws.Range("D43").Value = "67.39"
$ws.Range("E43").Value = "  -1.40%  "

$ws.Range("E44").Value = "  +6.65%  "

$ws.Range("D45").Value = "5.51"
$ws.Range("E45").Value = "  +2.83%  "

$ws.Range("E46").Value = "  +0.25%  "

$ws.Range("D47").Value = "1.784.12"
$ws.Range("E47").Value = "  +0.87%  "

$ws.Range("D48").Value = "87.80"
$ws.Range("E48").Value = "  +0.57%  "

$ws.Range("E49").Value = "  +0.81%  "

$ws.Range("E50").Value = "  +0.38%  "

$ws.Range("D51").Value = "7.60"
$ws.Range("E51").Value = "  -0.86%  "
